$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New rows 10-12: three extra vendor sub-tables (company, bank, purchasing)
# ---------------------------------------------------------------------

# Row 10 - S_SUPPL_COMPANY
$ws.Range("A10").Value = "vendor"
$ws.Range("B10").Value = "S_SUPPL_COMPANY"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = "LIFNR,BUKRS,AKONT,ZTERM1,ZWELS_01"

# Row 11 - S_SUPP_BANK
$ws.Range("A11").Value = "vendor"
$ws.Range("B11").Value = "S_SUPP_BANK"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = "LIFNR,BANKS,BANKL,BANKN,IBAN,BKONT,BKREF,KOINH,EBPP_ACCNAME"

# Row 12 - S_SUPPL_PURCHASING
$ws.Range("A12").Value = "vendor"
$ws.Range("B12").Value = "S_SUPPL_PURCHASING"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = "LIFNR,EKORG,WAERS,ZTERM,INCO1,INCO2,KALKS,VSBED,WEBRE,KZAUT,BSTAE,KZRET"

# ---------------------------------------------------------------------
# Formatting: ROWS / GEN_ORDER / COLUMN_ORDER columns stay numeric &
# right aligned, DOMAIN/TABLE_NAME stay general aligned (matches the
# rest of the table).
# ---------------------------------------------------------------------
$ws.Range("A10:B12").HorizontalAlignment = 1
$ws.Range("C10:D12").NumberFormat = "#,##0"
$ws.Range("C10:E12").HorizontalAlignment = -4152

# Row heights: the two new "sub" rows are slightly shorter than the
# standard data row, the final row matches the standard data row height.
$ws.Rows.Item(10).RowHeight = 17.25
$ws.Rows.Item(11).RowHeight = 17.25
$ws.Rows.Item(12).RowHeight = 19.5

# Widen column B so the longer structure names fit (target stored width
# 28.576428571428572; engine snaps ColumnWidth to 1/6-character pixel grid,
# so feed the nearest input that lands on that grid point).
$ws.Columns.Item(2).ColumnWidth = 27.666666666666668
